# Generate Report for Archive
# - Flip the "Ready for handoff" status to "In Translation" everywhere it
#   appears (Overview!E2/F2 summary columns, and the per-locale "Status"
#   column on the zh-cn / de-de sheets).
# - Shrink the now-narrower status columns to match the shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: zh-cn / de-de summary status columns (E, F)
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# Per-locale detail sheets: Status column (C)
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# Narrow the status columns to fit the shorter "In Translation" text.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
